# Master item sheet rework:
#  - Drop the separate numeric "id" counter column and the "key" column.
#  - The former "key" values (10000..10304) become the new "id" values in column A.
#  - Column B switches from the internal enum name (itemType) to a friendly "name".
#  - Column C (old "key") is removed entirely; data shifts left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New id (former key) values, in row order 2..10
$ids = @(10000, 10001, 10100, 10101, 10201, 10301, 10302, 10303, 10304)

# New friendly names, rows 2..9 (row 10's label is "Coin" and is left
# untouched below, since the text doesn't actually change there).
$names = @(
    "Item 1",
    "Item 2",
    "Puzzle Piece A",
    "Puzzle Piece B",
    "Old Key",
    "Book A",
    "Book B",
    "Book C"
)

$ws.Cells.Item(1, 8).Value = "itemType is an enum"

# Data rows: update the id (former key) column for every row, and the
# name column for every row except the last (Coin stays Coin).
for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

# Header row last: "name" is a brand-new label, so it's interned after all
# the other (re)used strings above.
$ws.Cells.Item(1, 1).Value = "id"
$ws.Cells.Item(1, 2).Value = "name"

# Clear the now-unused column C ("key") in place; columns D.. (and the H
# note column) keep their original positions/widths.
$ws.Range("C1:C10").Clear()

# Column A no longer needs the tight bestFit width from the old 1-digit ids;
# give it a fixed width to fit the new 5-digit ids.
$ws.Columns.Item(1).ColumnWidth = 10.21875

# Move the selection, matching the edited file's saved cursor position.
$ws.Range("D25").Select()
